# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I8").Value = 'sd'
$ws.Range("J8").Value = 'Statement-non-opinion'
$ws.Range("I9").Value = 'sd'
$ws.Range("J9").Value = 'Statement-non-opinion'
$ws.Range("I19").Value = 'ba'
$ws.Range("J19").Value = 'Appreciation'
$ws.Range("I21").Value = 'b'
$ws.Range("J21").Value = 'Acknowledge (Backchannel)'
$ws.Range("I22").Value = 'aa'
$ws.Range("J22").Value = 'Agree/Accept'
$ws.Range("I35").Value = 'aa'
$ws.Range("J35").Value = 'Agree/Accept'
$ws.Range("I37").Value = '%'
$ws.Range("J37").Value = 'Uninterpretable'
$ws.Range("I53").Value = 'sd'
$ws.Range("J53").Value = 'Statement-non-opinion'
$ws.Range("I54").Value = '%'
$ws.Range("J54").Value = 'Uninterpretable'
$ws.Range("I61").Value = 'sd'
$ws.Range("J61").Value = 'Statement-non-opinion'
$ws.Range("I66").Value = '%'
$ws.Range("J66").Value = 'Uninterpretable'
$ws.Range("I86").Value = 'sd'
$ws.Range("J86").Value = 'Statement-non-opinion'
$ws.Range("I111").Value = 'sd'
$ws.Range("J111").Value = 'Statement-non-opinion'
$ws.Range("I115").Value = 'aa'
$ws.Range("J115").Value = 'Agree/Accept'
$ws.Range("I116").Value = 'sd'
$ws.Range("J116").Value = 'Statement-non-opinion'
$ws.Range("I125").Value = 'ba'
$ws.Range("J125").Value = 'Appreciation'
$ws.Range("I130").Value = 'aa'
$ws.Range("J130").Value = 'Agree/Accept'
$ws.Range("I144").Value = 'aa'
$ws.Range("J144").Value = 'Agree/Accept'
$ws.Range("I155").Value = 'aa'
$ws.Range("J155").Value = 'Agree/Accept'
$ws.Range("I160").Value = 'b'
$ws.Range("J160").Value = 'Acknowledge (Backchannel)'
$ws.Range("I176").Value = 'b'
$ws.Range("J176").Value = 'Acknowledge (Backchannel)'
$ws.Range("I180").Value = 'sd'
$ws.Range("J180").Value = 'Statement-non-opinion'
$ws.Range("I183").Value = 'aa'
$ws.Range("J183").Value = 'Agree/Accept'
$ws.Range("I196").Value = 'b'
$ws.Range("J196").Value = 'Acknowledge (Backchannel)'
$ws.Range("I197").Value = 'ba'
$ws.Range("J197").Value = 'Appreciation'
$ws.Range("I249").Value = 'sv'
$ws.Range("J249").Value = 'Statement-opinion'
$ws.Range("I266").Value = '%'
$ws.Range("J266").Value = 'Uninterpretable'
$ws.Range("I275").Value = 'sd'
$ws.Range("J275").Value = 'Statement-non-opinion'
$ws.Range("I280").Value = 'aa'
$ws.Range("J280").Value = 'Agree/Accept'
$ws.Range("I295").Value = 'sd'
$ws.Range("J295").Value = 'Statement-non-opinion'
$ws.Range("I300").Value = 'sv'
$ws.Range("J300").Value = 'Statement-opinion'
$ws.Range("I302").Value = 'sv'
$ws.Range("J302").Value = 'Statement-opinion'
$ws.Range("I321").Value = 'aa'
$ws.Range("J321").Value = 'Agree/Accept'
$ws.Range("I328").Value = 'sv'
$ws.Range("J328").Value = 'Statement-opinion'
$ws.Range("I338").Value = 'sd'
$ws.Range("J338").Value = 'Statement-non-opinion'
$ws.Range("I343").Value = 'sd'
$ws.Range("J343").Value = 'Statement-non-opinion'
$ws.Range("I358").Value = 'b'
$ws.Range("J358").Value = 'Acknowledge (Backchannel)'
$ws.Range("I367").Value = 'aa'
$ws.Range("J367").Value = 'Agree/Accept'
$ws.Range("I392").Value = 'sv'
$ws.Range("J392").Value = 'Statement-opinion'
$ws.Range("I396").Value = 'sv'
$ws.Range("J396").Value = 'Statement-opinion'